$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while preserving it as TEXT
# (matching the original cell's stored type) instead of letting Excel
# auto-convert it to a number. We briefly mark the cell as Text ("@"),
# assign the value, then restore the original "Normal" style so the
# cell's formatting is left unchanged.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 13
Set-TextValue $ws.Range("B13") "34.13"
Set-TextValue $ws.Range("C13") "1.59"
Set-TextValue $ws.Range("D13") "35.71"

# Employment (% of total) - row 14
Set-TextValue $ws.Range("B14") "44.78"
Set-TextValue $ws.Range("C14") "36.46"
Set-TextValue $ws.Range("D14") "81.24"

# Enterprises (% of total) - row 16
Set-TextValue $ws.Range("B16") "95.44"
Set-TextValue $ws.Range("C16") "4.44"
Set-TextValue $ws.Range("D16") "99.88"
